$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update "Last Updated" timestamp ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "30 Oct 2025, 01:49 PM"

# --- "1 Month Performance" sheet: updated stock order / % change values ---
$wsPerf = $wb.Worksheets.Item("1 Month Performance")
$wsPerf.Range("C2").Value = 111.0513
$wsPerf.Range("C3").Value = 84.5455
$wsPerf.Range("C4").Value = 78.3373
$wsPerf.Range("C5").Value = 67.60809999999999
$wsPerf.Range("C8").Value = 53.8088
$wsPerf.Range("C9").Value = 50.6493
$wsPerf.Range("C10").Value = 45.0551
$wsPerf.Range("C11").Value = 42.1893
$wsPerf.Range("C14").Value = 39.0879
$wsPerf.Range("B15").Value = "SHAREINDIA"
$wsPerf.Range("C15").Value = 37.9488
$wsPerf.Range("B16").Value = "SANDUMA"
$wsPerf.Range("C16").Value = 37.9459
$wsPerf.Range("C17").Value = 36.3385
$wsPerf.Range("C18").Value = 35.485
$wsPerf.Range("C20").Value = 34.9934
$wsPerf.Range("C21").Value = 33.5418
$wsPerf.Range("C24").Value = 30.983
$wsPerf.Range("B25").Value = "ORIENTTECH"
$wsPerf.Range("C25").Value = 30.1699
$wsPerf.Range("B26").Value = "SOUTHBANK"
$wsPerf.Range("C26").Value = 30.1247
$wsPerf.Range("C28").Value = 28.4588
$wsPerf.Range("B29").Value = "MRPL"
$wsPerf.Range("C29").Value = 28.4117
$wsPerf.Range("B30").Value = "TARACHAND"
$wsPerf.Range("C30").Value = 28.3944
$wsPerf.Range("B31").Value = "TDPOWERSYS"
$wsPerf.Range("C31").Value = 27.236
$wsPerf.Range("B32").Value = "ADANIPOWER"
$wsPerf.Range("C32").Value = 27.1672
$wsPerf.Range("B33").Value = "EMKAY"
$wsPerf.Range("C33").Value = 27.0722
$wsPerf.Range("C34").Value = 26.195
$wsPerf.Range("B35").Value = "MARINE"
$wsPerf.Range("C35").Value = 25.4158
$wsPerf.Range("B36").Value = "HATSUN"
$wsPerf.Range("C36").Value = 25.4097
$wsPerf.Range("C37").Value = 25.2012
$wsPerf.Range("C40").Value = 24.609
$wsPerf.Range("C41").Value = 24.4179
$wsPerf.Range("C43").Value = 23.6754
$wsPerf.Range("C46").Value = 22.8035
$wsPerf.Range("C47").Value = 22.5444
$wsPerf.Range("C48").Value = 22.0093
$wsPerf.Range("C49").Value = 21.9672
$wsPerf.Range("C50").Value = 21.8294
$wsPerf.Range("B51").Value = "GRMOVER"
$wsPerf.Range("C51").Value = 21.5714
$wsPerf.Range("B52").Value = "SURYODAY"
$wsPerf.Range("C52").Value = 21.5706
$wsPerf.Range("B53").Value = "KERNEX"
$wsPerf.Range("C53").Value = 21.4868
$wsPerf.Range("B54").Value = "CPEDU"
$wsPerf.Range("C54").Value = 21.4026
$wsPerf.Range("B55").Value = "INDRAMEDCO"
$wsPerf.Range("C55").Value = 21.3868
$wsPerf.Range("B56").Value = "SCI"
$wsPerf.Range("C56").Value = 21.2684
$wsPerf.Range("B59").Value = "STYLAMIND"
$wsPerf.Range("C59").Value = 20.3924
$wsPerf.Range("B60").Value = "MOLDTECH"
$wsPerf.Range("C60").Value = 20.3828
$wsPerf.Range("B61").Value = "PRIVISCL"
$wsPerf.Range("C61").Value = 20.3325
$wsPerf.Range("C62").Value = 20.2766
$wsPerf.Range("C63").Value = 20.0358
$wsPerf.Range("C64").Value = 19.9862
$wsPerf.Range("B65").Value = "ETHOSLTD"
$wsPerf.Range("C65").Value = 19.835
$wsPerf.Range("B66").Value = "ASALCBR"
$wsPerf.Range("C66").Value = 19.825
$wsPerf.Range("B67").Value = "BHARATWIRE"
$wsPerf.Range("C67").Value = 19.698
$wsPerf.Range("B68").Value = "FEDERALBNK"
$wsPerf.Range("C68").Value = 19.6935
$wsPerf.Range("B69").Value = "BLUEDART"
$wsPerf.Range("C69").Value = 19.6598
$wsPerf.Range("C70").Value = 19.3982
$wsPerf.Range("B71").Value = "WHEELS"
$wsPerf.Range("C71").Value = 19.2837
$wsPerf.Range("B72").Value = "PRECWIRE"
$wsPerf.Range("C72").Value = 19.1473
$wsPerf.Range("C73").Value = 18.8776
$wsPerf.Range("C75").Value = 18.6973
$wsPerf.Range("C76").Value = 18.6947

# --- "distance from Dma50" sheet: updated distance values ---
$wsDma = $wb.Worksheets.Item("distance from Dma50")
$wsDma.Range("C2").Value = 9.6813
$wsDma.Range("C3").Value = 7.226
$wsDma.Range("C4").Value = 5.9191
$wsDma.Range("C5").Value = 5.0681
$wsDma.Range("C6").Value = 5.0565
$wsDma.Range("C8").Value = 4.4168
$wsDma.Range("C9").Value = 4.2681
$wsDma.Range("C10").Value = 3.7841
$wsDma.Range("C11").Value = 3.4249
$wsDma.Range("C12").Value = 3.2838
$wsDma.Range("C13").Value = 3.2688
$wsDma.Range("C14").Value = 2.9455
$wsDma.Range("C15").Value = 2.9099
$wsDma.Range("C16").Value = 2.8248
$wsDma.Range("C17").Value = 2.6829
$wsDma.Range("C18").Value = 2.4945
$wsDma.Range("C19").Value = 2.4465
$wsDma.Range("C20").Value = 2.1827
$wsDma.Range("C21").Value = 2.1523
$wsDma.Range("C22").Value = 1.3692
$wsDma.Range("C23").Value = 1.2592
$wsDma.Range("C24").Value = 1.2265
$wsDma.Range("C25").Value = 0.9661999999999999
$wsDma.Range("C26").Value = 0.9249000000000001
$wsDma.Range("C27").Value = 0.8846000000000001
$wsDma.Range("C28").Value = 0.5334
$wsDma.Range("C29").Value = 0.2078
$wsDma.Range("C30").Value = -2.1729
